$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - update "want to go" counts (column F) / min
# price (column G) for several existing events.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Cells.Item(2, 6).Value = 54
$wsExpo.Cells.Item(3, 6).Value = 276
$wsExpo.Cells.Item(3, 7).Value = 70
$wsExpo.Cells.Item(4, 6).Value = 969
$wsExpo.Cells.Item(7, 6).Value = 666
$wsExpo.Cells.Item(8, 6).Value = 235
$wsExpo.Cells.Item(10, 6).Value = 7
$wsExpo.Cells.Item(12, 6).Value = 179
$wsExpo.Cells.Item(13, 6).Value = 37
$wsExpo.Cells.Item(14, 6).Value = 767
$wsExpo.Cells.Item(15, 6).Value = 105
$wsExpo.Cells.Item(16, 6).Value = 1907
$wsExpo.Cells.Item(17, 6).Value = 422
$wsExpo.Cells.Item(18, 6).Value = 5783
$wsExpo.Cells.Item(19, 6).Value = 420
$wsExpo.Cells.Item(20, 6).Value = 509
$wsExpo.Cells.Item(21, 6).Value = 34
$wsExpo.Cells.Item(24, 6).Value = 183

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local Life) - update "want to go" counts (column F).
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item(3)
$wsLocal.Cells.Item(2, 6).Value = 5424
$wsLocal.Cells.Item(3, 6).Value = 361
$wsLocal.Cells.Item(4, 6).Value = 348

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All Types) - same counter updates as above, plus a brand
# new row (13) inserted for a newly scraped event, pushing every following
# row down by one.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Cells.Item(2, 6).Value = 54
$wsAll.Cells.Item(3, 6).Value = 5424
$wsAll.Cells.Item(4, 6).Value = 361
$wsAll.Cells.Item(6, 6).Value = 348
$wsAll.Cells.Item(7, 6).Value = 276
$wsAll.Cells.Item(7, 7).Value = 70

# Insert a new row at position 13; everything currently at/after row 13
# moves down to make room.
$wsAll.Rows.Item(13).Insert()

# Re-point the running index column (A) for every shifted row so that it
# keeps matching "row number - 1", then fill in the new row's data.
for ($r = 14; $r -le 49; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}

$a13 = $wsAll.Cells.Item(13, 1)
$a13.Value = 12
$a13.Font.Bold = $true
$a13.Borders.LineStyle = 1
$a13.HorizontalAlignment = -4108
$a13.VerticalAlignment = -4160

$b13 = $wsAll.Cells.Item(13, 2)
$b13.Value = "'2024-09-08"
$b13.Style = "Normal"

$wsAll.Cells.Item(13, 3).Value = "广州·神山羊2024巡演ENCOUNTER"
$wsAll.Cells.Item(13, 4).Value = "流花路117号流花展贸中心5号馆 广州大麦66live house"
$wsAll.Cells.Item(13, 5).Value = "2024.09.08 19:00-09.08 20:30"
$wsAll.Cells.Item(13, 6).Value = 495
$wsAll.Cells.Item(13, 7).Value = 380
$wsAll.Cells.Item(13, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89835"
$wsAll.Cells.Item(13, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/vpWr6GKQ1721877449091.jpeg"
